$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain plain text (they look like plain
# decimal numbers, so Excel would otherwise coerce them to numeric values).
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D41", "D43", "D45", "D46", "D48", "D49", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "26.222.25"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3
$ws.Range("D3").Value = "1.673.07"
$ws.Range("E3").Value = "  -1.34%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.72%  "

# Row 5
$ws.Range("D5").Value = "211.92"
$ws.Range("E5").Value = "  -2.98%  "

# Row 6
$ws.Range("D6").Value = "0.5282"
$ws.Range("E6").Value = "  -3.66%  "

# Row 8
$ws.Range("D8").Value = "0.2647"
$ws.Range("E8").Value = "  -3.20%  "

# Row 9
$ws.Range("D9").Value = "0.06287"
$ws.Range("E9").Value = "  -2.43%  "

# Row 10
$ws.Range("D10").Value = "21.34"
$ws.Range("E10").Value = "  -2.89%  "

# Row 11
$ws.Range("D11").Value = "0.07562"
$ws.Range("E11").Value = "  -1.49%  "

# Row 12
$ws.Range("D12").Value = "1.693.36"
$ws.Range("E12").Value = "  -0.05%  "

# Row 13
$ws.Range("D13").Value = "4.461"
$ws.Range("E13").Value = "  -2.10%  "

# Row 14
$ws.Range("D14").Value = "0.5601"
$ws.Range("E14").Value = "  -4.14%  "

# Row 15
$ws.Range("D15").Value = "67.08"
$ws.Range("E15").Value = "  +2.16%  "

# Row 16
$ws.Range("D16").Value = "0.000008029"

# Row 17
$ws.Range("D17").Value = "26.256.38"
$ws.Range("E17").Value = "  -0.80%  "

# Row 18
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19
$ws.Range("E19").Value = "  -2.89%  "

# Row 20
$ws.Range("D20").Value = "187.54"
$ws.Range("E20").Value = "  -1.88%  "

# Row 21
$ws.Range("D21").Value = "10.42"
$ws.Range("E21").Value = "  -5.22%  "

# Row 22
$ws.Range("D22").Value = "6.215"
$ws.Range("E22").Value = "  -0.71%  "

# Row 23
$ws.Range("E23").Value = "  -0.68%  "

# Row 24
$ws.Range("D24").Value = "149.97"
$ws.Range("E24").Value = "  +0.81%  "

# Row 25
$ws.Range("D25").Value = "0.1259"
$ws.Range("E25").Value = "  -4.16%  "

# Row 26
$ws.Range("D26").Value = "7.576"
$ws.Range("E26").Value = "  -4.32%  "

# Row 27
$ws.Range("D27").Value = "15.98"
$ws.Range("E27").Value = "  +1.17%  "

# Row 28
$ws.Range("D28").Value = "0.06199"
$ws.Range("E28").Value = "  -0.31%  "

# Row 29
$ws.Range("E29").Value = "  -1.73%  "

# Row 30
$ws.Range("E30").Value = "  -3.43%  "

# Row 31
$ws.Range("D31").Value = "3.505"
$ws.Range("E31").Value = "  -2.89%  "

# Row 32
$ws.Range("D32").Value = "3.432"
$ws.Range("E32").Value = "  -4.50%  "

# Row 33
$ws.Range("D33").Value = "1.635"
$ws.Range("E33").Value = "  -3.08%  "

# Row 34
$ws.Range("D34").Value = "1.003"
$ws.Range("E34").Value = "  -3.52%  "

# Row 35
$ws.Range("D35").Value = "0.6070"
$ws.Range("E35").Value = "  -1.34%  "

# Row 36
$ws.Range("E36").Value = "  -0.02%  "

# Row 37
$ws.Range("D37").Value = "2.746"
$ws.Range("E37").Value = "  -0.42%  "

# Row 38
$ws.Range("D38").Value = "6.120"
$ws.Range("E38").Value = "  +0.28%  "

# Row 39
$ws.Range("E39").Value = "  -1.90%  "

# Row 40
$ws.Range("D40").Value = "1.102.56"
$ws.Range("E40").Value = "  -1.36%  "

# Row 41
$ws.Range("D41").Value = "0.8758"
$ws.Range("E41").Value = "  -0.62%  "

# Row 43
$ws.Range("D43").Value = "99.97"
$ws.Range("E43").Value = "  -1.25%  "

# Row 44
$ws.Range("D44").Value = "1.823.54"

# Row 45
$ws.Range("D45").Value = "0.00000000108"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46
$ws.Range("D46").Value = "56.02"
$ws.Range("E46").Value = "  -2.80%  "

# Row 47
$ws.Range("E47").Value = "  -0.09%  "

# Row 48
$ws.Range("D48").Value = "8.005"
$ws.Range("E48").Value = "  -2.20%  "

# Row 49
$ws.Range("D49").Value = "0.05221"
$ws.Range("E49").Value = "  -1.17%  "

# Row 50
$ws.Range("E50").Value = "  -1.18%  "

# Row 51
$ws.Range("D51").Value = "5.992"
$ws.Range("E51").Value = "  -2.19%  "
